$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the date / volume / price values between row 2 and row 5
# Row 2 -> becomes what row 5 previously had
$ws.Range("D2").Value = 44316
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 1111

# Row 5 -> becomes what row 2 previously had
$ws.Range("D5").Value = 44516
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 33000
$ws.Range("O5").Value = 34000
$ws.Range("P5").Value = 33500
$ws.Range("S5").Value = 1861
